$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "310.39"
$ws.Range("E2").Value = "-3.31%"
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "54.42"
$ws.Range("E3").Value = "11.09%"
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.096"
$ws.Range("E4").Value = "-4.50%"
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.07912"
$ws.Range("E5").Value = "-1.71%"
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "4.547"
$ws.Range("E6").Value = "-0.95%"
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "1.398"
$ws.Range("E7").Value = "2.72%"
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "1.672"
$ws.Range("E8").Value = "2.13%"
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "0.1243"
$ws.Range("E9").Value = "-2.97%"
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.2020"
$ws.Range("E10").Value = "2.81%"
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.04742"
$ws.Range("E11").Value = "1.40%"
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.09437"
$ws.Range("E12").Value = "-2.30%"
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.1045"
$ws.Range("E13").Value = "0.01%"
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001271"
$ws.Range("E14").Value = "-3.34%"
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.005788"
$ws.Range("E15").Value = "-1.14%"
$rng.Style = "Normal"

$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$ws.Range("E16").Value = "2,016.70%"
$rng.Style = "Normal"

$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$ws.Range("E17").Value = "-0.61%"
$rng.Style = "Normal"

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = "-0.94%"
$rng.Style = "Normal"

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = "0.3428"
$ws.Range("E19").Value = "-2.23%"
$rng.Style = "Normal"

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = "8.408"
$ws.Range("E20").Value = "4.66%"
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "0.1359"
$ws.Range("E21").Value = "-0.38%"
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "0.2906"
$ws.Range("E22").Value = "-6.04%"
$rng.Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.04173"
$ws.Range("E23").Value = "-0.56%"
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.001253"
$ws.Range("E24").Value = "-4.71%"
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.003975"
$ws.Range("E25").Value = "-6.91%"
$rng.Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.0001345"
$ws.Range("E26").Value = "-0.43%"
$rng.Style = "Normal"

$rng = $ws.Range("D38:E38")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = "0.02628"
$ws.Range("E38").Value = "-3.54%"
$rng.Style = "Normal"

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.05929"
$ws.Range("E39").Value = "-2.40%"
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.01081"
$ws.Range("E40").Value = "-0.52%"
$rng.Style = "Normal"

$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$ws.Range("E41").Value = "19.53%"
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.007928"
$ws.Range("E42").Value = "-1.16%"
$rng.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.008155"
$ws.Range("E43").Value = "3.21%"
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.008323"
$ws.Range("E44").Value = "-3.96%"
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.3385"
$ws.Range("E45").Value = "-3.10%"
$rng.Style = "Normal"

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00007257"
$ws.Range("E46").Value = "5.62%"
$rng.Style = "Normal"

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000747"
$ws.Range("E47").Value = "-0.36%"
$rng.Style = "Normal"

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.06007"
$ws.Range("E48").Value = "2.23%"
$rng.Style = "Normal"

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = "0.002609"
$ws.Range("E49").Value = "-34.75%"
$rng.Style = "Normal"

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002091"
$ws.Range("E50").Value = "-0.36%"
$rng.Style = "Normal"

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = "0.0001991"
$ws.Range("E51").Value = "-0.36%"
$rng.Style = "Normal"
